# Add a new contact row (row 6) to the Contacts sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / numeric fields can be assigned directly.
$ws.Cells.Item(6, 1).Value = "Андрей Андреев"      # FullName
$ws.Cells.Item(6, 2).Value = 28                     # Age
$ws.Cells.Item(6, 3).Value = 3                       # Experience
$ws.Cells.Item(6, 4).Value = "andreyandreev@mail.ru" # Email

# The phone number starts with "+" and is otherwise all digits, so a plain
# assignment would be auto-converted to a number (losing the "+"). Build it
# as a text formula in a scratch cell, then paste just the resulting value
# (not the formula) into the target cell so it keeps its text type without
# altering the cell's number format / style.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="+380957894561"'
$scratch.Copy()
$ws.Cells.Item(6, 5).PasteSpecial(-4163) # xlPasteValues
$scratch.ClearContents()
$excel.CutCopyMode = 0
